$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the "Undertaken a Xamarin forms ... UDemy ..." bullet paragraph by
# scanning the bullet list (searching by text keeps this resilient to the
# exact paragraph index).
# ---------------------------------------------------------------------------
$xamarinIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Undertaken a")) {
        $xamarinIndex = $i
        break
    }
}
if ($xamarinIndex -eq -1) {
    throw "Could not find the 'Undertaken a Xamarin forms' bullet paragraph"
}

$xamarinPara = $d.Paragraphs.Item($xamarinIndex)

# ---------------------------------------------------------------------------
# Step 1: insert a brand-new bullet paragraph right after it - this will
# become the relocated "Undertaken a Xamarin forms ... UDemy ..." bullet.
# ---------------------------------------------------------------------------
$xamarinPara.Range.InsertParagraphAfter()

$newParaIndex = $xamarinIndex + 1
$newPara = $d.Paragraphs.Item($newParaIndex)
$newRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$newRange.Text = "Undertaken a Xamarin forms (.Net standard, PCL) course in my own time via UDemy, developed a small app to try out new skills."

# Newly created paragraphs/ranges in this host default their character
# formatting to bold, so explicitly clear it before re-applying bold only
# where it belongs.
$newPara = $d.Paragraphs.Item($newParaIndex)
$newWhole = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$newWhole.Font.Bold = $False

$newPara = $d.Paragraphs.Item($newParaIndex)
$newFull = $newPara.Range.Text
$newStart = $newPara.Range.Start
$boldText = "Xamarin forms (.Net standard, PCL)"
$boldStart = $newStart + $newFull.IndexOf($boldText)
$boldEnd = $boldStart + $boldText.Length
$boldRange = $d.Range($boldStart, $boldEnd)
$boldRange.Font.Bold = $True

# ---------------------------------------------------------------------------
# Step 2: replace the original bullet's content with the new JavaScript
# frameworks sentence ("Experienced with a few JavaScript frameworks
# (React, knockout, vanilla ...)"), with the framework list in bold.
# ---------------------------------------------------------------------------
$xamarinPara = $d.Paragraphs.Item($xamarinIndex)
$origRange = $d.Range($xamarinPara.Range.Start, $xamarinPara.Range.End - 1)
$origRange.Text = "Experienced with a few JavaScript frameworks (React, knockout, vanilla `u{2026})"

$xamarinPara = $d.Paragraphs.Item($xamarinIndex)
$origWhole = $d.Range($xamarinPara.Range.Start, $xamarinPara.Range.End - 1)
$origWhole.Font.Bold = $False

$xamarinPara = $d.Paragraphs.Item($xamarinIndex)
$jsFull = $xamarinPara.Range.Text
$jsStart = $xamarinPara.Range.Start
$jsBoldText = "React, knockout, vanilla `u{2026})"
$jsBoldStart = $jsStart + $jsFull.IndexOf($jsBoldText)
$jsBoldEnd = $jsBoldStart + $jsBoldText.Length
$jsBoldRange = $d.Range($jsBoldStart, $jsBoldEnd)
$jsBoldRange.Font.Bold = $True
